# Generate Report for Handoff
# The file "98ddb7bb-b688-4c50-8d50-1ea5bc4d443b.md" (row 3 on every sheet) has
# just been handed off for localization, so its status moves from
# "In Translation" to "Ready for handoff", its Priority moves from "ht" to
# "mt", and the various "Latest Handoff"/"Latest HO Xliff Generate" timestamps
# are refreshed.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Overview sheet: row 3 is the 98ddb7bb...md file.
#   E3 = zh-cn status, F3 = de-de status, G3 = Latest HO Xliff Generate Date
# ---------------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E3").Value = "Ready for handoff"
$overview.Range("F3").Value = "Ready for handoff"
$overview.Range("G3").Value = "2016-09-06 10:18:28"

# Status column widened to fit the longer "Ready for handoff" text.
$overview.Columns.Item(5).ColumnWidth = 16.33
$overview.Columns.Item(6).ColumnWidth = 16.33

# ---------------------------------------------------------------------------
# zh-cn sheet: row 3 is the 98ddb7bb...md file.
#   C3 = Status, E3 = Priority, H3 = Latest Handoff Datetime
# ---------------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = "Ready for handoff"
$zhcn.Range("E3").Value = "mt"
$zhcn.Range("H3").Value = "2016-09-06 10:18:24"

$zhcn.Columns.Item(3).ColumnWidth = 16.33

# ---------------------------------------------------------------------------
# de-de sheet: row 3 is the 98ddb7bb...md file.
#   C3 = Status, E3 = Priority, H3 = Latest Handoff Datetime
# ---------------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = "Ready for handoff"
$dede.Range("E3").Value = "mt"
$dede.Range("H3").Value = "2016-09-06 10:18:28"

$dede.Columns.Item(3).ColumnWidth = 16.33
